$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A9 with corrected precise value
$ws.Range("A9").Value = 45864.5419865162

# Add new row 10 with sensor data
$ws.Range("A10").Value = 45864.58361671767
$ws.Range("A10").NumberFormat = $ws.Range("A9").NumberFormat

$ws.Range("B10").Value = 2025
$ws.Range("C10").Value = 30
$ws.Range("D10").Value = 21.44
$ws.Range("E10").Value = 67.34
$ws.Range("F10").Value = 101.95
$ws.Range("G10").Value = 16.32
$ws.Range("H10").Value = "SE"
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = "14:00:24"
